# Apply metadata updates to the "Metadata" worksheet of the CodeSystem workbook.
# Changes:
#   URL:       http://ibm.com/...               -> http://linuxforhealth.org/...
#   Version:   7.0.0                            -> 8.0.0
#   Date:      2022-09-08T16:11:15+00:00        -> 2022-11-10T16:00:46+00:00
#   Publisher: Alvearie Team                    -> LinuxForHealth Team

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

$ws.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/CodeSystem/eng-communication-request-status-reason"
$ws.Range("B3").Value = "8.0.0"
$ws.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$ws.Range("B9").Value = "LinuxForHealth Team"
